$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row fixes: H1 drops leading space; I1 gets new "Autosave" header ---
$ws.Range("H1").Value = "Success"
$ws.Range("I1").Value = "Autosave"

# --- Append new experiment rows 71-86 ---
# (Leading apostrophe on True/False values forces literal text, matching
#  the source data which stores them as inline strings, not booleans.)
$ws.Range("A71").Value = "sin(10*x)"
$ws.Range("B71").Value = 20000
$ws.Range("C71").Value = 4
$ws.Range("D71").Value = 4
$ws.Range("E71").Value = 0.0007102209492586553
$ws.Range("F71").Value = 81.3016402721405
$ws.Range("G71").Value = 245.9974968900273
$ws.Range("H71").Value = "'True"
$ws.Range("I71").Value = "'False"
$ws.Range("A72").Value = "sin(10*x)"
$ws.Range("B72").Value = 20000
$ws.Range("C72").Value = 4
$ws.Range("D72").Value = 4
$ws.Range("E72").Value = 0.0007102209492586553
$ws.Range("F72").Value = 77.89687538146973
$ws.Range("G72").Value = 256.749707893388
$ws.Range("H72").Value = "'True"
$ws.Range("I72").Value = "'False"
$ws.Range("A73").Value = "sin(10*x)"
$ws.Range("B73").Value = 20000
$ws.Range("C73").Value = 4
$ws.Range("D73").Value = 4
$ws.Range("E73").Value = 0.0007102209492586553
$ws.Range("F73").Value = 82.5798716545105
$ws.Range("G73").Value = 242.1897685149478
$ws.Range("H73").Value = "'True"
$ws.Range("I73").Value = "'False"
$ws.Range("A74").Value = "25*x**5 - 50*x**4 + 100*x**3 - 200*x**2 + 400*x - 800"
$ws.Range("B74").Value = 20000
$ws.Range("C74").Value = 4
$ws.Range("D74").Value = 4
$ws.Range("E74").Value = 22.01958465576172
$ws.Range("F74").Value = 74.77325892448425
$ws.Range("G74").Value = 267.4753018348257
$ws.Range("H74").Value = "'True"
$ws.Range("I74").Value = "'True"
$ws.Range("A75").Value = "5*x**5 - 5*x**4 + 5*x**3 - 5*x**2 + 5*x - 5"
$ws.Range("B75").Value = 20000
$ws.Range("C75").Value = 4
$ws.Range("D75").Value = 4
$ws.Range("E75").Value = 0.4573033154010773
$ws.Range("F75").Value = 77.77278709411621
$ws.Range("G75").Value = 257.1593580129401
$ws.Range("H75").Value = "'True"
$ws.Range("I75").Value = "'True"
$ws.Range("A76").Value = "5*x**5 - 5*x**4 + 5*x**3 - 5*x**2 + 5*x - 5"
$ws.Range("B76").Value = 20000
$ws.Range("C76").Value = 4
$ws.Range("D76").Value = 4
$ws.Range("E76").Value = 0.04460600018501282
$ws.Range("F76").Value = 74.65597915649414
$ws.Range("G76").Value = 267.8954884253266
$ws.Range("H76").Value = "'True"
$ws.Range("I76").Value = "'True"
$ws.Range("A77").Value = "5*x**5 - 5*x**4 + 5*x**3 - 5*x**2 + 5*x - 5"
$ws.Range("B77").Value = 5000
$ws.Range("C77").Value = 4
$ws.Range("D77").Value = 4
$ws.Range("E77").Value = 0.04698513075709343
$ws.Range("F77").Value = 19.63630437850952
$ws.Range("G77").Value = 254.630398043337
$ws.Range("H77").Value = "'True"
$ws.Range("I77").Value = "'True"
$ws.Range("A78").Value = "5*x**5 - 5*x**4 + 5*x**3 - 5*x**2 + 5*x - 5"
$ws.Range("B78").Value = 5000
$ws.Range("C78").Value = 4
$ws.Range("D78").Value = 4
$ws.Range("E78").Value = 0.04698513075709343
$ws.Range("F78").Value = 20.19875931739807
$ws.Range("G78").Value = 247.5399563622347
$ws.Range("H78").Value = "'True"
$ws.Range("I78").Value = "'True"
$ws.Range("A79").Value = "5*x**5 - 5*x**4 + 5*x**3 - 5*x**2 + 5*x - 5"
$ws.Range("B79").Value = 5000
$ws.Range("C79").Value = 4
$ws.Range("D79").Value = 4
$ws.Range("E79").Value = 0.04698513075709343
$ws.Range("F79").Value = 20.41918420791626
$ws.Range("G79").Value = 244.8677649943313
$ws.Range("H79").Value = "'True"
$ws.Range("I79").Value = "'True"
$ws.Range("A80").Value = "sin(10*x)"
$ws.Range("B80").Value = 5000
$ws.Range("C80").Value = 1
$ws.Range("D80").Value = 1
$ws.Range("E80").Value = 0.0001324334734817967
$ws.Range("F80").Value = 8.728749513626099
$ws.Range("G80").Value = 572.8197369158895
$ws.Range("H80").Value = "'True"
$ws.Range("I80").Value = "'True"
$ws.Range("A81").Value = "tan(10*x)"
$ws.Range("B81").Value = 5000
$ws.Range("C81").Value = 1
$ws.Range("D81").Value = 1
$ws.Range("E81").Value = 3846.375244140625
$ws.Range("F81").Value = 8.680053234100342
$ws.Range("G81").Value = 576.0333335695531
$ws.Range("H81").Value = "'True"
$ws.Range("I81").Value = "'True"
$ws.Range("A82").Value = "tan(x)"
$ws.Range("B82").Value = 5000
$ws.Range("C82").Value = 1
$ws.Range("D82").Value = 1
$ws.Range("E82").Value = 0.00001114379301725421
$ws.Range("F82").Value = 8.594872951507568
$ws.Range("G82").Value = 581.7421651500949
$ws.Range("H82").Value = "'True"
$ws.Range("I82").Value = "'True"
$ws.Range("A83").Value = "tan(1.57*x)"
$ws.Range("B83").Value = 5000
$ws.Range("C83").Value = 1
$ws.Range("D83").Value = 1
$ws.Range("E83").Value = 3147.503173828125
$ws.Range("F83").Value = 8.639046907424927
$ws.Range("G83").Value = 578.7675485015243
$ws.Range("H83").Value = "'True"
$ws.Range("I83").Value = "'True"
$ws.Range("A84").Value = "tan(1.57*x)"
$ws.Range("B84").Value = 5000
$ws.Range("C84").Value = 4
$ws.Range("D84").Value = 4
$ws.Range("E84").Value = 2884.274169921875
$ws.Range("F84").Value = 20.30728888511658
$ws.Range("G84").Value = 246.2170124375663
$ws.Range("H84").Value = "'True"
$ws.Range("I84").Value = "'True"
$ws.Range("A85").Value = "sin(10*x)"
$ws.Range("B85").Value = 5000
$ws.Range("C85").Value = 1
$ws.Range("D85").Value = 1
$ws.Range("E85").Value = 0.0001320796355219922
$ws.Range("F85").Value = 17.51730537414551
$ws.Range("G85").Value = 285.4320281120234
$ws.Range("H85").Value = "'True"
$ws.Range("I85").Value = "'True"
$ws.Range("A86").Value = "sin(10*x)"
$ws.Range("B86").Value = 5000
$ws.Range("C86").Value = 1
$ws.Range("D86").Value = 1
$ws.Range("E86").Value = 0.0001320943410973996
$ws.Range("F86").Value = 18.81007742881775
$ws.Range("G86").Value = 265.8149610984488
$ws.Range("H86").Value = "'True"
$ws.Range("I86").Value = "'False"
